# Scheduled runner refresh: update cached Universalis price-check columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# eight crafting-job sheets with freshly pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 836.7143
$ws.Range("I9").Value = 891.4
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 891.4
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = -722.4
$ws.Range("N9").Value = -1038
$ws.Range("H55").Value = 126.666664
$ws.Range("J55").Value = 90
$ws.Range("L55").Value = 90
$ws.Range("N55").Value = -518
$ws.Range("H98").Value = 4869.7144
$ws.Range("I98").Value = 1651.8636
$ws.Range("K98").Value = 1651.8636
$ws.Range("M98").Value = -153.8635999999999
$ws.Range("H100").Value = 1089.25
$ws.Range("I100").Value = 1152
$ws.Range("K100").Value = 1152
$ws.Range("M100").Value = -611
$ws.Range("H112").Value = 6114.2354
$ws.Range("I112").Value = 1865.3334
$ws.Range("J112").Value = 6525.4194
$ws.Range("K112").Value = 5596.0002
$ws.Range("L112").Value = 19576.2582
$ws.Range("M112").Value = -4488.0002
$ws.Range("N112").Value = -21792.2582
$ws.Range("H113").Value = 7091.615
$ws.Range("I113").Value = 7362.25
$ws.Range("J113").Value = 6658.6
$ws.Range("K113").Value = 7362.25
$ws.Range("L113").Value = 6658.6
$ws.Range("M113").Value = -4108.25
$ws.Range("N113").Value = -13166.6
$ws.Range("H122").Value = 4869.7144
$ws.Range("I122").Value = 1651.8636
$ws.Range("K122").Value = 4955.5908
$ws.Range("M122").Value = -2505.5908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 19666.666
$ws.Range("I33").Value = 19666.666
$ws.Range("K33").Value = 19666.666
$ws.Range("M33").Value = -19337.666
$ws.Range("H61").Value = 6477.304
$ws.Range("J61").Value = 28650
$ws.Range("L61").Value = 28650
$ws.Range("N61").Value = -29074
$ws.Range("H74").Value = 360578.47
$ws.Range("I74").Value = 1200522
$ws.Range("J74").Value = 10602
$ws.Range("K74").Value = 1200522
$ws.Range("L74").Value = 10602
$ws.Range("M74").Value = -1199648
$ws.Range("N74").Value = -12350
$ws.Range("H77").Value = 360578.47
$ws.Range("I77").Value = 1200522
$ws.Range("J77").Value = 10602
$ws.Range("K77").Value = 6002610
$ws.Range("L77").Value = 53010
$ws.Range("M77").Value = -5998242
$ws.Range("N77").Value = -61746
$ws.Range("H109").Value = 70999.5
$ws.Range("J109").Value = 70999.5
$ws.Range("L109").Value = 70999.5
$ws.Range("N109").Value = -73773.5
$ws.Range("H122").Value = 1960.6666
$ws.Range("I122").Value = 1685.8334
$ws.Range("J122").Value = 2785.1667
$ws.Range("K122").Value = 5057.5002
$ws.Range("L122").Value = 8355.500100000001
$ws.Range("M122").Value = -2607.5002
$ws.Range("N122").Value = -13255.5001
$ws.Range("H136").Value = 6477.304
$ws.Range("J136").Value = 28650
$ws.Range("L136").Value = 85950
$ws.Range("N136").Value = -91050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 878.875
$ws.Range("I80").Value = 924.75
$ws.Range("J80").Value = 833
$ws.Range("K80").Value = 924.75
$ws.Range("L80").Value = 833
$ws.Range("M80").Value = 73.25
$ws.Range("N80").Value = -2829
$ws.Range("H83").Value = 878.875
$ws.Range("I83").Value = 924.75
$ws.Range("J83").Value = 833
$ws.Range("K83").Value = 4623.75
$ws.Range("L83").Value = 4165
$ws.Range("M83").Value = 368.25
$ws.Range("N83").Value = -14149

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4249.5
$ws.Range("I3").Value = 4249.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4249.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4136.5
$ws.Range("N3").ClearContents()
$ws.Range("H16").Value = 1377.4814
$ws.Range("I16").Value = 1290.4
$ws.Range("J16").Value = 1626.2858
$ws.Range("K16").Value = 1290.4
$ws.Range("L16").Value = 1626.2858
$ws.Range("M16").Value = -1003.4
$ws.Range("N16").Value = -2200.2858
$ws.Range("H31").Value = 7693836.5
$ws.Range("J31").Value = 1787.6
$ws.Range("L31").Value = 1787.6
$ws.Range("N31").Value = -2377.6
$ws.Range("H34").Value = 7693836.5
$ws.Range("J34").Value = 1787.6
$ws.Range("L34").Value = 1787.6
$ws.Range("N34").Value = -2191.6
$ws.Range("H94").Value = 1626.174
$ws.Range("J94").Value = 2154.9092
$ws.Range("L94").Value = 2154.9092
$ws.Range("N94").Value = -3056.9092
$ws.Range("H99").Value = 5927.5
$ws.Range("I99").Value = 4622.3
$ws.Range("K99").Value = 4622.3
$ws.Range("M99").Value = -3124.3
$ws.Range("H105").Value = 1648.1666
$ws.Range("I105").Value = 977.9
$ws.Range("K105").Value = 977.9
$ws.Range("M105").Value = 769.1
$ws.Range("H113").Value = 1377.4814
$ws.Range("I113").Value = 1290.4
$ws.Range("J113").Value = 1626.2858
$ws.Range("K113").Value = 1290.4
$ws.Range("L113").Value = 1626.2858
$ws.Range("M113").Value = 879.5999999999999
$ws.Range("N113").Value = -5966.2858
$ws.Range("H126").Value = 5927.5
$ws.Range("I126").Value = 4622.3
$ws.Range("K126").Value = 13866.9
$ws.Range("M126").Value = -11396.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 869.6667
$ws.Range("J5").Value = 1099
$ws.Range("L5").Value = 3297
$ws.Range("N5").Value = -3521
$ws.Range("H14").Value = 153.18182
$ws.Range("I14").Value = 153.18182
$ws.Range("K14").Value = 459.5454599999999
$ws.Range("M14").Value = -286.5454599999999
$ws.Range("H23").Value = 436.46155
$ws.Range("I23").Value = 264
$ws.Range("J23").Value = 584.2857
$ws.Range("K23").Value = 792
$ws.Range("L23").Value = 1752.8571
$ws.Range("M23").Value = -557
$ws.Range("N23").Value = -2222.8571
$ws.Range("H40").Value = 42.142857
$ws.Range("I40").Value = 24.166666
$ws.Range("K40").Value = 96.666664
$ws.Range("M40").Value = -27.666664
$ws.Range("H121").Value = 98327.14
$ws.Range("I121").Value = 133528.75
$ws.Range("J121").Value = 51391.668
$ws.Range("K121").Value = 400586.25
$ws.Range("L121").Value = 154175.004
$ws.Range("M121").Value = -399276.25
$ws.Range("N121").Value = -156795.004
$ws.Range("H127").Value = 4795.25
$ws.Range("J127").Value = 5333.143
$ws.Range("L127").Value = 15999.429
$ws.Range("N127").Value = -25919.429
$ws.Range("H135").Value = 869.6667
$ws.Range("J135").Value = 1099
$ws.Range("L135").Value = 9891
$ws.Range("N135").Value = -14961

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7187.375
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 8916.5
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 8916.5
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -10912.5
$ws.Range("H83").Value = 7187.375
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 8916.5
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 44582.5
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -54566.5
$ws.Range("H113").Value = 1345.0667
$ws.Range("J113").Value = 1566.3334
$ws.Range("L113").Value = 1566.3334
$ws.Range("N113").Value = -5906.3334
$ws.Range("H122").Value = 4178.9375
$ws.Range("I122").Value = 3872.2693
$ws.Range("J122").Value = 5507.8335
$ws.Range("K122").Value = 11616.8079
$ws.Range("L122").Value = 16523.5005
$ws.Range("M122").Value = -9166.8079
$ws.Range("N122").Value = -21423.5005
$ws.Range("H132").Value = 2566.4285
$ws.Range("J132").Value = 2965.8572
$ws.Range("L132").Value = 8897.571599999999
$ws.Range("N132").Value = -13957.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 13750
$ws.Range("J43").Value = 13750
$ws.Range("L43").Value = 13750
$ws.Range("N43").Value = -14136
$ws.Range("H55").Value = 2423.9285
$ws.Range("I55").Value = 1667.8572
$ws.Range("J55").Value = 3180
$ws.Range("K55").Value = 1667.8572
$ws.Range("L55").Value = 3180
$ws.Range("M55").Value = -1494.8572
$ws.Range("N55").Value = -3526
$ws.Range("H68").Value = 2569.8333
$ws.Range("I68").Value = 2569.8333
$ws.Range("K68").Value = 2569.8333
$ws.Range("M68").Value = -1820.8333
$ws.Range("H71").Value = 2569.8333
$ws.Range("I71").Value = 2569.8333
$ws.Range("K71").Value = 12849.1665
$ws.Range("M71").Value = -9105.166499999999
$ws.Range("H93").Value = 945.3200000000001
$ws.Range("I93").Value = 946.5
$ws.Range("J93").Value = 940.6
$ws.Range("K93").Value = 946.5
$ws.Range("L93").Value = 940.6
$ws.Range("M93").Value = 301.5
$ws.Range("N93").Value = -3436.6
$ws.Range("H122").Value = 3500.1904
$ws.Range("I122").Value = 3412.5881
$ws.Range("K122").Value = 10237.7643
$ws.Range("M122").Value = -7787.764299999999
$ws.Range("H132").Value = 4230.75
$ws.Range("I132").Value = 3899.4614
$ws.Range("K132").Value = 11698.3842
$ws.Range("M132").Value = -9168.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10648.167
$ws.Range("I3").Value = 7963
$ws.Range("K3").Value = 7963
$ws.Range("M3").Value = -7849
$ws.Range("H30").Value = 20000
$ws.Range("I30").Value = 20000
$ws.Range("K30").Value = 20000
$ws.Range("M30").Value = -19893
$ws.Range("H74").Value = 18238.6
$ws.Range("I74").Value = 10569
$ws.Range("J74").Value = 20156
$ws.Range("K74").Value = 10569
$ws.Range("L74").Value = 20156
$ws.Range("M74").Value = -9633
$ws.Range("N74").Value = -22028
$ws.Range("H77").Value = 18238.6
$ws.Range("I77").Value = 10569
$ws.Range("J77").Value = 20156
$ws.Range("K77").Value = 31707
$ws.Range("L77").Value = 60468
$ws.Range("M77").Value = -27027
$ws.Range("N77").Value = -69828
